$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - F column holds "想去人数" (want-to-go count)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 0
$wsExhibit.Range("F6").Value = 0
$wsExhibit.Range("F8").Value = 0
$wsExhibit.Range("F10").Value = 0

# Sheet "全部类型" (All types) - F column holds "想去人数" (want-to-go count)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 0
$wsAll.Range("F3").Value = 116
$wsAll.Range("F4").Value = 1629
$wsAll.Range("F5").Value = 0
$wsAll.Range("F6").Value = 0
$wsAll.Range("F8").Value = 0
$wsAll.Range("F10").Value = 492

$wb.Save()
